$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text (avoids Excel auto-converting
# number-looking strings and losing trailing zeros / exact formatting)
function Set-Text($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# Row 2 - Bitcoin
Set-Text "D2" "66.262.04"
$ws.Range("E2").Value = "  +1.37%  "

# Row 3 - Ethereum
Set-Text "D3" "3.562.17"
$ws.Range("E3").Value = "  +4.59%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-Text "D5" "607.16"
$ws.Range("E5").Value = "  +2.31%  "

# Row 6 - Solana
Set-Text "D6" "145.50"
$ws.Range("E6").Value = "  +2.31%  "

# Row 7 - LidoStakedEther
Set-Text "D7" "3.560.44"
$ws.Range("E7").Value = "  +4.67%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.16%  "

# Row 9 - XRP
Set-Text "D9" "0.485"
$ws.Range("E9").Value = "  +3.99%  "

# Row 10 - Dogecoin
Set-Text "D10" "0.137"
$ws.Range("E10").Value = "  +1.60%  "

# Row 11 - Toncoin
Set-Text "D11" "7.99"
$ws.Range("E11").Value = "  +1.59%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.77%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-Text "D13" "4.165.80"

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +4.57%  "

# Row 15 - Avalanche
Set-Text "D15" "30.26"
$ws.Range("E15").Value = "  +1.56%  "

# Row 16 - WrappedEther
Set-Text "D16" "3.563.12"
$ws.Range("E16").Value = "  +4.61%  "

# Row 17 - WrappedBTC
Set-Text "D17" "66.356.34"
$ws.Range("E17").Value = "  +1.34%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.67%  "

# Row 19 - Uniswap
Set-Text "D19" "11.48"
$ws.Range("E19").Value = "  +10.98%  "

# Row 20 - Polkadot
Set-Text "D20" "6.23"
$ws.Range("E20").Value = "  +2.01%  "

# Row 21 - Chainlink
Set-Text "D21" "15.00"
$ws.Range("E21").Value = "  +1.31%  "

# Row 22 - BitcoinCash
Set-Text "D22" "431.63"
$ws.Range("E22").Value = "  +3.52%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +4.85%  "

# Row 24 - Litecoin
Set-Text "D24" "78.60"
$ws.Range("E24").Value = "  +1.96%  "

# Row 25 - WrappedeETH
Set-Text "D25" "3.705.53"
$ws.Range("E25").Value = "  +4.58%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +8.27%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  +4.72%  "

# Row 29 - RenderToken
Set-Text "D29" "8.07"
$ws.Range("E29").Value = "  +3.29%  "

# Row 30 - InternetComputer(DFINITY)
Set-Text "D30" "9.23"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  -0.12%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +1.89%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  -0.99%  "

# Row 34 - RenzoRestakedETH
Set-Text "D34" "3.558.59"
$ws.Range("E34").Value = "  +4.50%  "

# Row 35 - EthereumClassic
Set-Text "D35" "25.43"
$ws.Range("E35").Value = "  +3.39%  "

# Rows 36 & 37 - swap USDe <-> ImmutableX (ImmutableX now ranked above USDe)
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-Text "D36" "1.76"
$ws.Range("E36").Value = "  +3.09%  "

$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-Text "D37" "1.00"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38 - Aptos
$ws.Range("E38").Value = "  +4.48%  "

# Row 39 - NEARProtocol
$ws.Range("E39").Value = "  +2.47%  "

# Row 40 - FirstDigitalUSD
Set-Text "D40" "0.999"
$ws.Range("E40").Value = "  -0.03%  "

# Row 41 - Monero
Set-Text "D41" "170.93"
$ws.Range("E41").Value = "  -1.35%  "

# Row 42 - Hedera
$ws.Range("E42").Value = "  -0.48%  "

# Row 43 - Filecoin
Set-Text "D43" "5.23"
$ws.Range("E43").Value = "  +3.65%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  +3.48%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +1.77%  "

# Row 46 - OKB
Set-Text "D46" "46.27"
$ws.Range("E46").Value = "  +1.63%  "

# Rows 47 & 48 - swap InjectiveProtocol <-> ONDO (ONDO now ranked above InjectiveProtocol)
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-Text "D47" "1.22"
$ws.Range("E47").Value = "  +4.72%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-Text "D48" "26.16"
$ws.Range("E48").Value = "  -1.90%  "

# Row 49 - dogwifhat
$ws.Range("E49").Value = "  +5.73%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +0.79%  "

# Row 51 - SuiNetwork
Set-Text "D51" "0.957"
$ws.Range("E51").Value = "  +4.10%  "
